# "Fruta / hortaliza, semanal" — the upstream weekly re-pull reshuffled
# the 5 Cereza price rows (A2:T6) for Agrícola del Norte S.A. de Arica.
# Re-seat each row's A:T values at its new position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of each data row exactly as it stands today, keyed by its
# current row number.
$data = @{
    2 = @{ A=1; B='Agrícola del Norte S.A. de Arica'; C='Arica y Parinacota'; D=44229; E=15; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Santina'; L='Primera'; M=250; N=6500;  O=7000;  P=6750;  Q='$/bandeja 5 kilos';  R='Provincia de Curicó';    S=1350; T=5  }
    3 = @{ A=1; B='Agrícola del Norte S.A. de Arica'; C='Arica y Parinacota'; D=44161; E=15; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Bing';    L='Primera'; M=160; N=39000; O=40000; P=39500; Q='$/caja 20 kilos';   R='Provincia de Curicó';    S=1975; T=20 }
    4 = @{ A=1; B='Agrícola del Norte S.A. de Arica'; C='Arica y Parinacota'; D=44208; E=15; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Lapins';  L='Segunda'; M=200; N=10500; O=11000; P=10750; Q='$/bandeja 12 kilos'; R='Provincia de Curicó';    S=896;  T=12 }
    5 = @{ A=1; B='Agrícola del Norte S.A. de Arica'; C='Arica y Parinacota'; D=44210; E=15; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Rainier'; L='Segunda'; M=250; N=21000; O=22000; P=21500; Q='$/caja 18 kilos';   R="Región de O'Higgins"; S=1194; T=18 }
    6 = @{ A=1; B='Agrícola del Norte S.A. de Arica'; C='Arica y Parinacota'; D=44175; E=15; F='Fruta'; G=100103; H='Frutos de hueso (carozo)'; I=100103001; J='Cereza'; K='Rainier'; L='Segunda'; M=270; N=25000; O=26000; P=25500; Q='$/caja 18 kilos';   R="Región de O'Higgins"; S=1417; T=18 }
}

# New row <- old row, per the refreshed feed.
$newRowFromOld = @{ 2 = 6; 3 = 5; 4 = 2; 5 = 3; 6 = 4 }

$cols = @('A','B','C','D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T')

foreach ($oldRow in $newRowFromOld.Keys) {
    $newRow = $newRowFromOld[$oldRow]
    $row = $data[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $row[$col]
    }
}
